$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 and row 4 each get a new Young's-modulus entry in column G and a
# new (zero) entry in column L, matching what's already entered in row 2
# for the same two columns (E = 3.0e11, and a 0 in the new L column).

$ws.Range("G3").Value = 300000000000
$ws.Range("G4").Value = 300000000000
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 0

# Copy the existing number formatting from row 2's G/L cells onto the
# freshly entered cells so they render the same way (G uses the 2-decimal
# numeric style already used throughout that column; L matches the style
# already used by L2).
$ws.Range("G2").Copy()
$ws.Range("G3:G4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("L2").Copy()
$ws.Range("L3:L4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# The saved cursor/selection moved from N3 to L4.
$ws.Range("L4").Select()
